$wb = $excel.ActiveWorkbook

# Update the "Date" metadata value (Metadata sheet, B8)
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2025-10-29T11:46:56+00:00"

# Fix casing / wording of the ExerciceProfessionnel references on the
# "Elements" sheet (row 6 -> the Specialite.exerciceProfessionnel element)
$wsElem = $wb.Worksheets.Item("Elements")
$wsElem.Range("A6").Value = "Specialite.ExerciceProfessionnel"
$wsElem.Range("B6").Value = "Specialite.ExerciceProfessionnel"
$wsElem.Range("L6").Value = "Lien vers la classe ExerciceProfessionnel"
$wsElem.Range("M6").Value = "Lien vers la classe ExerciceProfessionnel"
$wsElem.Range("AF6").Value = "SavoirFaire.ExerciceProfessionnel"

# Slight bestFit column-width adjustments (Elements sheet) caused by the
# text-length change above (26.2265625 -> 26.25390625 and 27.3984375 ->
# 27.42578125 in the saved OOXML "width" units)
$wsElem.Columns.Item(1).ColumnWidth = 25.5
$wsElem.Columns.Item(2).ColumnWidth = 25.5
$wsElem.Columns.Item(32).ColumnWidth = 26.7
